# Insert a new weekly data row at row 253 (pushing the existing rows
# 253..294 down to 254..295) on the single data sheet of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 253:294 down by one row, duplicating formatting from the
# row above (Excel's default Insert behaviour) so the date cell in the
# new row keeps the same custom date number format as its neighbours.
$ws.Rows.Item(253).Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A253").Value = 6
$ws.Range("B253").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C253").Value = "Metropolitana"
$ws.Range("D253").Value = 44543
$ws.Range("E253").Value = 13
$ws.Range("F253").Value = 100112032
$ws.Range("G253").Value = "Zapallo italiano"
$ws.Range("H253").Value = "Sin especificar"
$ws.Range("I253").Value = "Primera"
$ws.Range("J253").Value = 230
$ws.Range("K253").Value = 8000
$ws.Range("L253").Value = 9000
$ws.Range("M253").Value = 8565
$ws.Range("N253").Value = "`$/caja 50 unidades"
$ws.Range("O253").Value = "Región de O'Higgins"
$ws.Range("P253").Value = 171
$ws.Range("Q253").Value = 50
$ws.Range("R253").Value = "Hortaliza"
